$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 98
$ws.Range("F98").Value = 13886961
$ws.Range("G98").Value = 11076397

# Row 99
$ws.Range("M99").Value = 8586181

# Row 105
$ws.Range("M105").Value = 8358621

# Row 132
$ws.Range("H132").Value = 935347
$ws.Range("M132").Value = 14066894

# Row 133
$ws.Range("M133").Value = 14469347

# Row 134
$ws.Range("M134").Value = 14604982

# Row 135
$ws.Range("M135").Value = 14690419

# Row 136
$ws.Range("M136").Value = 14866878

# Row 137
$ws.Range("M137").Value = 14861396

# Row 138
$ws.Range("M138").Value = 14670788

# Row 139
$ws.Range("M139").Value = 14467013

# Row 140
$ws.Range("D140").Value = 823721
$ws.Range("M140").Value = 14131793
$ws.Range("N140").Value = 8935557

# Row 141
$ws.Range("M141").Value = 13894089
$ws.Range("N141").Value = 8782072

# Row 142
$ws.Range("M142").Value = 13749996
$ws.Range("N142").Value = 8776340

# Row 143
$ws.Range("M143").Value = 13621865
$ws.Range("N143").Value = 8836782

# Row 144
$ws.Range("N144").Value = 9042590

# Row 145
$ws.Range("N145").Value = 9395210

# Row 146
$ws.Range("N146").Value = 9716074

# Row 147
$ws.Range("N147").Value = 9916905

# Row 148
$ws.Range("N148").Value = 10176538

# Row 149
$ws.Range("N149").Value = 10402119

# Row 150
$ws.Range("N150").Value = 10586519

# Row 151
$ws.Range("N151").Value = 10834432

# Row 156
$ws.Range("J156").Value = 37.93

# Row 186
$ws.Range("D186").Value = 590378
$ws.Range("G186").Value = 9816802
$ws.Range("H186").Value = 1209084
$ws.Range("I186").Value = 127.39
$ws.Range("M186").Value = 13889874
$ws.Range("N186").Value = 8211183

# Row 187
$ws.Range("M187").Value = 13966106
$ws.Range("N187").Value = 8231174

# Row 188
$ws.Range("L188").Value = 125.9
$ws.Range("M188").Value = 14177340
$ws.Range("N188").Value = 8310960

# Row 189
$ws.Range("M189").Value = 14260648
$ws.Range("N189").Value = 8315320

# Row 190
$ws.Range("M190").Value = 14422975
$ws.Range("N190").Value = 8370901

# Row 191
$ws.Range("D191").Value = 573544
$ws.Range("G191").Value = 8670388
$ws.Range("H191").Value = 1073051
$ws.Range("M191").Value = 14528123
$ws.Range("N191").Value = 8368419

# Row 192
$ws.Range("M192").Value = 14532898
$ws.Range("N192").Value = 8289993

# Row 193
$ws.Range("M193").Value = 14761421
$ws.Range("N193").Value = 8365085

# Row 194
$ws.Range("F194").Value = 14729052
$ws.Range("M194").Value = 14844102
$ws.Range("N194").Value = 8387972

# Row 195
$ws.Range("M195").Value = 14889033
$ws.Range("N195").Value = 8325444

# Row 196
$ws.Range("F196").Value = 13097783
$ws.Range("M196").Value = 14967714
$ws.Range("N196").Value = 8319940

# Row 197
$ws.Range("M197").Value = 15089859
$ws.Range("N197").Value = 8341404

# Row 198
$ws.Range("M198").Value = 15141500
$ws.Range("N198").Value = 8330304

# Row 199
$ws.Range("M199").Value = 15309043
$ws.Range("N199").Value = 8356552

# Row 200
$ws.Range("M200").Value = 15440653
$ws.Range("N200").Value = 8350571

# Row 201
$ws.Range("E201").Value = 210443
$ws.Range("M201").Value = 15510428
$ws.Range("N201").Value = 8368314
$ws.Range("O201").Value = 2605563

# Row 202
$ws.Range("F202").Value = 10071208
$ws.Range("M202").Value = 15673738
$ws.Range("N202").Value = 8381871
$ws.Range("O202").Value = 2635476

# Row 203
$ws.Range("F203").Value = 9851482
$ws.Range("O203").Value = 2650975

# Row 204
$ws.Range("O204").Value = 2658797

# Row 205
$ws.Range("O205").Value = 2687346

# Row 206
$ws.Range("O206").Value = 2682400

# Row 207
$ws.Range("O207").Value = 2694751

# Row 208
$ws.Range("O208").Value = 2719334

# Row 209
$ws.Range("O209").Value = 2707776

# Row 210
$ws.Range("O210").Value = 2680341

# Row 211
$ws.Range("O211").Value = 2671979

# Row 212
$ws.Range("O212").Value = 2637236

# Row 232
$ws.Range("B232").Value = 1361137
$ws.Range("F232").Value = 11589025

# Row 624
$ws.Range("B624").Value = 1151677
$ws.Range("D624").Value = 328617
$ws.Range("E624").Value = 211223
$ws.Range("F624").Value = 7849467
$ws.Range("H624").Value = 2223070
$ws.Range("M624").Value = 48468526
$ws.Range("N624").Value = 8090235
$ws.Range("O624").Value = 5602474

# Row 625
$ws.Range("B625").Value = 1249893
$ws.Range("D625").Value = 413887
$ws.Range("E625").Value = 182251
$ws.Range("F625").Value = 7672045
$ws.Range("G625").Value = 6097705
$ws.Range("H625").Value = 2192414
$ws.Range("I625").Value = 376.43
$ws.Range("K625").Value = 17.81
$ws.Range("L625").Value = 349.81
$ws.Range("M625").Value = 44271472
$ws.Range("N625").Value = 7629053
$ws.Range("O625").Value = 4743274

# Row 626
$ws.Range("H626").Value = 2634737
$ws.Range("L626").Value = 352.7
$ws.Range("M626").Value = 41287840
$ws.Range("N626").Value = 7246641
$ws.Range("O626").Value = 4099379

# Row 627
$ws.Range("H627").Value = 2415448
$ws.Range("L627").Value = 355.84
$ws.Range("M627").Value = 38635131
$ws.Range("N627").Value = 6872121
$ws.Range("O627").Value = 3605431
$ws.Range("R627").Value = 324761

# Row 628
$ws.Range("G628").Value = 6440467
$ws.Range("H628").Value = 2541590
$ws.Range("L628").Value = 359.75
$ws.Range("M628").Value = 35830260
$ws.Range("N628").Value = 6391881
$ws.Range("O628").Value = 3165221
$ws.Range("P628").Value = 5213
$ws.Range("Q628").Value = 14087
$ws.Range("R628").Value = 325304

# Row 629
$ws.Range("B629").Value = 877536
$ws.Range("C629").Value = 1430
$ws.Range("D629").Value = 304514
$ws.Range("E629").Value = 125435
$ws.Range("F629").Value = 6045274
$ws.Range("G629").Value = 4792493
$ws.Range("H629").Value = 1887964
$ws.Range("I629").Value = 404
$ws.Range("J629").Value = 33.89
$ws.Range("K629").Value = 17.41
$ws.Range("L629").Value = 363.4
$ws.Range("M629").Value = 33649585
$ws.Range("N629").Value = 5630122
$ws.Range("O629").Value = 2881924
$ws.Range("P629").Value = 5224
$ws.Range("Q629").Value = 14108
$ws.Range("R629").Value = 323956

# Row 630
$ws.Range("A630").Value = 44712
$ws.Range("B630").Value = 859766
$ws.Range("D630").Value = 312842
$ws.Range("E630").Value = 125662
$ws.Range("F630").Value = 5997528
$ws.Range("G630").Value = 4677140
$ws.Range("H630").Value = 1815690
$ws.Range("I630").Value = 397.16
$ws.Range("J630").Value = 32.6
$ws.Range("K630").Value = 17.42
$ws.Range("L630").Value = 367.86
$ws.Range("M630").Value = 31547709
$ws.Range("N630").Value = 5190810
$ws.Range("O630").Value = 2637041
$ws.Range("P630").Value = 5260
$ws.Range("Q630").Value = 14124
